$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - worksheet 1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3313
$ws1.Range("F6").Value = 2270
$ws1.Range("F12").Value = 29
$ws1.Range("F14").Value = 331

# Sheet "全部类型" (All types) - worksheet 4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3313
$ws4.Range("F7").Value = 2270
$ws4.Range("F13").Value = 29
$ws4.Range("F15").Value = 331
